$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before the summary-formula block (originally rows 18-20)
# so the sheet grows from 17 data rows (2-17) to 25 data rows (2-26); the
# COUNTIF/SUM formulas shift down from rows 18-20 to rows 27-29 automatically.
for ($i = 0; $i -lt 9; $i++) {
    $ws.Rows.Item(18).Insert()
}

# Column F holds textual percentages ("82%"); force text formatting first so
# Excel does not auto-convert the assigned strings into numeric percentages.
$ws.Range("F2:F26").NumberFormat = "@"

# Row 2: Morocco ✓ - Republic of the Congo: 1:0
$ws.Cells.Item(2,1).Value = "Tue Oct 14"
$ws.Cells.Item(2,2).Value = "Morocco ✓ - Republic of the Congo: 1:0"
$ws.Cells.Item(2,3).Value = 1.87
$ws.Cells.Item(2,4).Value = "Morocco"
$ws.Cells.Item(2,5).Value = 2.5
$ws.Cells.Item(2,6).Value = "82%"
$ws.Cells.Item(2,7).Value = "✓"
$ws.Cells.Item(2,8).Value = 1
$ws.Cells.Item(2,9).Value = $true

# Row 3: Latvia - England ✓: 0:5
$ws.Cells.Item(3,1).Value = "Tue Oct 14"
$ws.Cells.Item(3,2).Value = "Latvia - England ✓: 0:5"
$ws.Cells.Item(3,3).Value = 1.39
$ws.Cells.Item(3,4).Value = "England"
$ws.Cells.Item(3,5).Value = 2.5
$ws.Cells.Item(3,6).Value = "80%"
$ws.Cells.Item(3,7).Value = "✓"
$ws.Cells.Item(3,8).Value = 5
$ws.Cells.Item(3,9).Value = $false

# Row 4: Norway  - New Zealand: 1:1
$ws.Cells.Item(4,1).Value = "Tue Oct 14"
$ws.Cells.Item(4,2).Value = "Norway  - New Zealand: 1:1"
$ws.Cells.Item(4,3).Value = 2.9
$ws.Cells.Item(4,4).Value = "Norway"
$ws.Cells.Item(4,5).Value = 3.5
$ws.Cells.Item(4,6).Value = "79%"
$ws.Cells.Item(4,7).ClearContents()
$ws.Cells.Item(4,8).Value = 2
$ws.Cells.Item(4,9).Value = $true

# Row 5: Russia ✓ - Bolivia: 3:0
$ws.Cells.Item(5,1).Value = "Tue Oct 14"
$ws.Cells.Item(5,2).Value = "Russia ✓ - Bolivia: 3:0"
$ws.Cells.Item(5,3).Value = 4.25
$ws.Cells.Item(5,4).Value = "Russia"
$ws.Cells.Item(5,5).Value = 5.5
$ws.Cells.Item(5,6).Value = "79%"
$ws.Cells.Item(5,7).Value = "✓"
$ws.Cells.Item(5,8).Value = 3
$ws.Cells.Item(5,9).Value = $true

# Row 6: Puerto Rico - Argentina : 00:00
$ws.Cells.Item(6,1).Value = "Tue Oct 14"
$ws.Cells.Item(6,2).Value = "Puerto Rico - Argentina : 00:00"
$ws.Cells.Item(6,3).Value = 1.66
$ws.Cells.Item(6,4).Value = "Argentina"
$ws.Cells.Item(6,5).Value = 2.5
$ws.Cells.Item(6,6).Value = "79%"
$ws.Cells.Item(6,7).ClearContents()
$ws.Cells.Item(6,8).Value = 0
$ws.Cells.Item(6,9).Value = $true

# Row 7: Spain ✓ - Bulgaria: 4:0
$ws.Cells.Item(7,1).Value = "Tue Oct 14"
$ws.Cells.Item(7,2).Value = "Spain ✓ - Bulgaria: 4:0"
$ws.Cells.Item(7,3).Value = 6.2
$ws.Cells.Item(7,4).Value = "Spain"
$ws.Cells.Item(7,5).Value = 7.5
$ws.Cells.Item(7,6).Value = "79%"
$ws.Cells.Item(7,7).Value = "✓"
$ws.Cells.Item(7,8).Value = 4
$ws.Cells.Item(7,9).Value = $true

# Row 8: Portugal  - Hungary: 2:2
$ws.Cells.Item(8,1).Value = "Tue Oct 14"
$ws.Cells.Item(8,2).Value = "Portugal  - Hungary: 2:2"
$ws.Cells.Item(8,3).Value = 3.47
$ws.Cells.Item(8,4).Value = "Portugal"
$ws.Cells.Item(8,5).Value = 4.5
$ws.Cells.Item(8,6).Value = "77%"
$ws.Cells.Item(8,7).ClearContents()
$ws.Cells.Item(8,8).Value = 4
$ws.Cells.Item(8,9).Value = $true

# Row 9: Senegal ✓ - Mauritania: 4:0
$ws.Cells.Item(9,1).Value = "Tue Oct 14"
$ws.Cells.Item(9,2).Value = "Senegal ✓ - Mauritania: 4:0"
$ws.Cells.Item(9,3).Value = 0.2
$ws.Cells.Item(9,4).Value = "Senegal"
$ws.Cells.Item(9,5).Value = 1.5
$ws.Cells.Item(9,6).Value = "74%"
$ws.Cells.Item(9,7).Value = "✓"
$ws.Cells.Item(9,8).Value = 4
$ws.Cells.Item(9,9).Value = $false

# Row 10: Italy ✓ - Israel: 3:0
$ws.Cells.Item(10,1).Value = "Tue Oct 14"
$ws.Cells.Item(10,2).Value = "Italy ✓ - Israel: 3:0"
$ws.Cells.Item(10,3).Value = 4.09
$ws.Cells.Item(10,4).Value = "Italy"
$ws.Cells.Item(10,5).Value = 5.5
$ws.Cells.Item(10,6).Value = "74%"
$ws.Cells.Item(10,7).Value = "✓"
$ws.Cells.Item(10,8).Value = 3
$ws.Cells.Item(10,9).Value = $true

# Row 11: Malaysia ✓ - Laos: 5:1
$ws.Cells.Item(11,1).Value = "Tue Oct 14"
$ws.Cells.Item(11,2).Value = "Malaysia ✓ - Laos: 5:1"
$ws.Cells.Item(11,3).Value = 2.54
$ws.Cells.Item(11,4).Value = "Malaysia"
$ws.Cells.Item(11,5).Value = 3.5
$ws.Cells.Item(11,6).Value = "73%"
$ws.Cells.Item(11,7).Value = "✓"
$ws.Cells.Item(11,8).Value = 6
$ws.Cells.Item(11,9).Value = $false

# Row 12: Democratic Republic of the Congo ✓ - Sudan: 1:0
$ws.Cells.Item(12,1).Value = "Tue Oct 14"
$ws.Cells.Item(12,2).Value = "Democratic Republic of the Congo ✓ - Sudan: 1:0"
$ws.Cells.Item(12,3).Value = 0.61
$ws.Cells.Item(12,4).Value = "Democratic Republic of the Congo"
$ws.Cells.Item(12,5).Value = 1.5
$ws.Cells.Item(12,6).Value = "73%"
$ws.Cells.Item(12,7).Value = "✓"
$ws.Cells.Item(12,8).Value = 1
$ws.Cells.Item(12,9).Value = $true

# Row 13: Ivory Coast ✓ - Kenya: 3:0
$ws.Cells.Item(13,1).Value = "Tue Oct 14"
$ws.Cells.Item(13,2).Value = "Ivory Coast ✓ - Kenya: 3:0"
$ws.Cells.Item(13,3).Value = 0.9
$ws.Cells.Item(13,4).Value = "Ivory Coast"
$ws.Cells.Item(13,5).Value = 1.5
$ws.Cells.Item(13,6).Value = "73%"
$ws.Cells.Item(13,7).Value = "✓"
$ws.Cells.Item(13,8).Value = 3
$ws.Cells.Item(13,9).Value = $false

# Row 14: Gabon ✓ - Burundi: 2:0
$ws.Cells.Item(14,1).Value = "Tue Oct 14"
$ws.Cells.Item(14,2).Value = "Gabon ✓ - Burundi: 2:0"
$ws.Cells.Item(14,3).Value = 0.81
$ws.Cells.Item(14,4).Value = "Gabon"
$ws.Cells.Item(14,5).Value = 1.5
$ws.Cells.Item(14,6).Value = "73%"
$ws.Cells.Item(14,7).Value = "✓"
$ws.Cells.Item(14,8).Value = 2
$ws.Cells.Item(14,9).Value = $false

# Row 15: Algeria ✓ - Uganda: 2:1
$ws.Cells.Item(15,1).Value = "Tue Oct 14"
$ws.Cells.Item(15,2).Value = "Algeria ✓ - Uganda: 2:1"
$ws.Cells.Item(15,3).Value = 1.34
$ws.Cells.Item(15,4).Value = "Algeria"
$ws.Cells.Item(15,5).Value = 2.5
$ws.Cells.Item(15,6).Value = "72%"
$ws.Cells.Item(15,7).Value = "✓"
$ws.Cells.Item(15,8).Value = 3
$ws.Cells.Item(15,9).Value = $false

# Row 16: Nigeria ✓ - Benin: 4:0
$ws.Cells.Item(16,1).Value = "Tue Oct 14"
$ws.Cells.Item(16,2).Value = "Nigeria ✓ - Benin: 4:0"
$ws.Cells.Item(16,3).Value = 0.92
$ws.Cells.Item(16,4).Value = "Nigeria"
$ws.Cells.Item(16,5).Value = 1.5
$ws.Cells.Item(16,6).Value = "71%"
$ws.Cells.Item(16,7).Value = "✓"
$ws.Cells.Item(16,8).Value = 4
$ws.Cells.Item(16,9).Value = $false

# Row 17: FBC Melgar  - Alianza Universidad: 00:00
$ws.Cells.Item(17,1).Value = "Tue Oct 14"
$ws.Cells.Item(17,2).Value = "FBC Melgar  - Alianza Universidad: 00:00"
$ws.Cells.Item(17,3).Value = 2.25
$ws.Cells.Item(17,4).Value = "FBC Melgar"
$ws.Cells.Item(17,5).Value = 3.5
$ws.Cells.Item(17,6).Value = "70%"
$ws.Cells.Item(17,7).ClearContents()
$ws.Cells.Item(17,8).Value = 0
$ws.Cells.Item(17,9).Value = $true

# Row 18: Philippines ✓ - Timor-Leste: 3:1
$ws.Cells.Item(18,1).Value = "Tue Oct 14"
$ws.Cells.Item(18,2).Value = "Philippines ✓ - Timor-Leste: 3:1"
$ws.Cells.Item(18,3).Value = 3.02
$ws.Cells.Item(18,4).Value = "Philippines"
$ws.Cells.Item(18,5).Value = 4.5
$ws.Cells.Item(18,6).Value = "70%"
$ws.Cells.Item(18,7).Value = "✓"
$ws.Cells.Item(18,8).Value = 4
$ws.Cells.Item(18,9).Value = $true

# Row 19: Nepal - Vietnam ✓: 0:1
$ws.Cells.Item(19,1).Value = "Tue Oct 14"
$ws.Cells.Item(19,2).Value = "Nepal - Vietnam ✓: 0:1"
$ws.Cells.Item(19,3).Value = 1.5
$ws.Cells.Item(19,4).Value = "Vietnam"
$ws.Cells.Item(19,5).Value = 2.5
$ws.Cells.Item(19,6).Value = "68%"
$ws.Cells.Item(19,7).Value = "✓"
$ws.Cells.Item(19,8).Value = 1
$ws.Cells.Item(19,9).Value = $true

# Row 20: Guinea  - Botswana: 2:2
$ws.Cells.Item(20,1).Value = "Tue Oct 14"
$ws.Cells.Item(20,2).Value = "Guinea  - Botswana: 2:2"
$ws.Cells.Item(20,3).Value = 0.49
$ws.Cells.Item(20,4).Value = "Guinea"
$ws.Cells.Item(20,5).Value = 1.5
$ws.Cells.Item(20,6).Value = "68%"
$ws.Cells.Item(20,7).ClearContents()
$ws.Cells.Item(20,8).Value = 4
$ws.Cells.Item(20,9).Value = $false

# Row 21: Costa Rica ✓ - Nicaragua: 4:1
$ws.Cells.Item(21,1).Value = "Tue Oct 14"
$ws.Cells.Item(21,2).Value = "Costa Rica ✓ - Nicaragua: 4:1"
$ws.Cells.Item(21,3).Value = 2.31
$ws.Cells.Item(21,4).Value = "Costa Rica"
$ws.Cells.Item(21,5).Value = 3.5
$ws.Cells.Item(21,6).Value = "67%"
$ws.Cells.Item(21,7).Value = "✓"
$ws.Cells.Item(21,8).Value = 5
$ws.Cells.Item(21,9).Value = $false

# Row 22: Iran ✓ - Tanzania: 2:0
$ws.Cells.Item(22,1).Value = "Tue Oct 14"
$ws.Cells.Item(22,2).Value = "Iran ✓ - Tanzania: 2:0"
$ws.Cells.Item(22,3).Value = 0.46
$ws.Cells.Item(22,4).Value = "Iran"
$ws.Cells.Item(22,5).Value = 1.5
$ws.Cells.Item(22,6).Value = "66%"
$ws.Cells.Item(22,7).Value = "✓"
$ws.Cells.Item(22,8).Value = 2
$ws.Cells.Item(22,9).Value = $false

# Row 23: South Africa ✓ - Rwanda: 3:0
$ws.Cells.Item(23,1).Value = "Tue Oct 14"
$ws.Cells.Item(23,2).Value = "South Africa ✓ - Rwanda: 3:0"
$ws.Cells.Item(23,3).Value = 0.25
$ws.Cells.Item(23,4).Value = "South Africa"
$ws.Cells.Item(23,5).Value = 1.5
$ws.Cells.Item(23,6).Value = "61%"
$ws.Cells.Item(23,7).Value = "✓"
$ws.Cells.Item(23,8).Value = 3
$ws.Cells.Item(23,9).Value = $false

# Row 24: Albania ✓ - Jordan: 4:2
$ws.Cells.Item(24,1).Value = "Tue Oct 14"
$ws.Cells.Item(24,2).Value = "Albania ✓ - Jordan: 4:2"
$ws.Cells.Item(24,3).Value = 0
$ws.Cells.Item(24,4).Value = "Albania"
$ws.Cells.Item(24,5).Value = 0.5
$ws.Cells.Item(24,6).Value = "59%"
$ws.Cells.Item(24,7).Value = "✓"
$ws.Cells.Item(24,8).Value = 6
$ws.Cells.Item(24,9).Value = $false

# Row 25: Chinese Taipei - Thailand ✓: 1:6
$ws.Cells.Item(25,1).Value = "Tue Oct 14"
$ws.Cells.Item(25,2).Value = "Chinese Taipei - Thailand ✓: 1:6"
$ws.Cells.Item(25,3).Value = 3.23
$ws.Cells.Item(25,4).Value = "Thailand"
$ws.Cells.Item(25,5).Value = 4.5
$ws.Cells.Item(25,6).Value = "58%"
$ws.Cells.Item(25,7).Value = "✓"
$ws.Cells.Item(25,8).Value = 7
$ws.Cells.Item(25,9).Value = $false

# Row 26: Turkiye ✓ - Georgia: 4:1
$ws.Cells.Item(26,1).Value = "Tue Oct 14"
$ws.Cells.Item(26,2).Value = "Turkiye ✓ - Georgia: 4:1"
$ws.Cells.Item(26,3).Value = 2.46
$ws.Cells.Item(26,4).Value = "Turkiye"
$ws.Cells.Item(26,5).Value = 3.5
$ws.Cells.Item(26,6).Value = "57%"
$ws.Cells.Item(26,7).Value = "✓"
$ws.Cells.Item(26,8).Value = 5
$ws.Cells.Item(26,9).Value = $false
